# This workbook tracks daily Cilantro price records. A new daily record
# needs to be inserted as row 191, pushing all subsequent records down by
# one row (the previously-last row 312 becomes row 313).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 191 (shifts rows 191..312 down to 192..313)
$ws.Rows.Item(191).Insert()

# Populate the newly-inserted row 191 with the new daily record
$ws.Cells.Item(191, 1).Value  = 10
$ws.Cells.Item(191, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(191, 3).Value  = "La Araucanía"
$ws.Cells.Item(191, 4).Value  = 44606
$ws.Cells.Item(191, 5).Value  = 9
$ws.Cells.Item(191, 6).Value  = 100112040
$ws.Cells.Item(191, 7).Value  = "Cilantro"
$ws.Cells.Item(191, 8).Value  = "Sin especificar"
$ws.Cells.Item(191, 9).Value  = "Primera"
$ws.Cells.Item(191, 10).Value = 220
$ws.Cells.Item(191, 11).Value = 5000
$ws.Cells.Item(191, 12).Value = 6000
$ws.Cells.Item(191, 13).Value = 5545
$ws.Cells.Item(191, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(191, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(191, 16).Value = 2772
$ws.Cells.Item(191, 17).Value = 2
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Match the style of the date cell (D) in the neighboring row so the new
# row's date cell keeps the same date number format as the rest of the column.
$ws.Cells.Item(191, 4).NumberFormat = $ws.Cells.Item(192, 4).NumberFormat

Write-Host "Inserted new record row at 191; dimension should now be A1:R313"
